$d = $word.ActiveDocument

# Locate the paragraph that ends with the Jira-tickets sentence.
$r = $d.Content
$found = $r.Find.Execute(
    "Submitted and completed over 1000 Jira tickets during my tenure at Merchology.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the anchor paragraph text."
}

# Collapse the found range to its end and insert a brand-new paragraph
# mark right after it. The new paragraph automatically inherits the
# same paragraph formatting (the bullet list: numPr ilvl=0 numId=1001)
# as the paragraph it was split off from.
$r.Collapse(0)
$r.InsertParagraphAfter()

# The newly created paragraph is empty (its Range.Text is just the
# paragraph mark, so Length is 1). Find it by scanning the document's
# paragraphs collection.
$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Length -eq 1) {
        $newPara = $p
    }
}

if (-not $newPara) {
    throw "Could not locate the newly inserted paragraph."
}

$newPara.Range.Text = "Chaired the Giveback Committee, an employee-managed group that coordinated and scheduled volunteer activities for Merchology employees. Collectively we volunteered over 1000 hours in 2023."

Write-Output "OK: inserted new bullet paragraph after the Jira-tickets item."
